$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width in the source XML is 16.42578125 character-units. The COM
# ColumnWidth setter here quantizes stored width to multiples of 1/6, so the
# nearest achievable value (16.5) is used (15.666666666666666 -> stored 16.5).
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666

$values = @(
    0.10989254549611616,
    -0.0059999999245725633,
    -0.0039999999377364759,
    -0.0079999998848041542,
    -0.002999999943356535,
    -0.0019999999486657316,
    -0.0099999998513333743,
    -0.0099999998516429045,
    -0.0019999999542061886,
    -0.0019999999582349659,
    -0.0029999999461356452,
    -0.0034999999412690386,
    -0.0034999999491720502,
    -0.0079999998969304542,
    -0.00099999998752053898,
    -0.0019999999804767299,
    0.036633801316864556,
    -0.0039999999668793862,
    0.0028114148340452338,
    -0.0005892022106817052,
    0.0091833430406094152,
    -0.0039999999494808591,
    -0.0049999999211944868,
    -0.019999999728490536,
    -0.019999999724596762,
    -0.002499999941893094,
    -0.0024999999404680118,
    -0.001999999941160624,
    -0.0069999998757790394,
    -0.059999999219685574,
    -0.0069999998853074175,
    0.055644490481812525,
    -0.003999999925984099
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

Write-Output "applied"
